$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "10000000=Red de 08 Nodos"

# Delete the "Fuerza Brutal" rows (row 8 in phi-section, row 16 in ms-section)
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(8).Delete()

$ws.Range("B1").Value = "ABCDEFGH=(ABCDEFGH|ABCDEFGH)"

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0

$ws.Range("B9").Value = 0.04808282852172852
$ws.Range("B10").Value = 1260.479565143585
$ws.Range("B11").Value = 166.4057974815369
$ws.Range("B12").Value = 69.00500822067261
$ws.Range("B13").Value = 25.3839955329895
$ws.Range("B14").Value = 25.20000004768372
